$wb = $excel.ActiveWorkbook

# --- meta_data sheet: selection D1 -> D59 ---
$wsMeta = $wb.Worksheets.Item("meta_data")
$wsMeta.Range("D59").Select()

# --- bioenergetics_control sheet: value change + tabSelected removed + selection B8 -> D14 ---
$wsBio = $wb.Worksheets.Item("bioenergetics_control")
$wsBio.Range("B6").Value = 0.083500000000000005
$wsBio.Range("D14").Select()

# --- Pyrs sheet: delete row 2 (shifts all years/data up by one row) + selection R81 -> D2 ---
$wsPyrs = $wb.Worksheets.Item("Pyrs")
$wsPyrs.Rows(2).Delete()
$wsPyrs.Range("D2").Select()

# --- UobsWtAge sheet: topLeftCell cleared + selection R142 -> I2 ---
$wsUobsWt = $wb.Worksheets.Item("UobsWtAge")
$wsUobsWt.Range("I2").Select()

# --- control sheet: value change + becomes the active/selected tab + selection B2 -> B21 ---
$wsControl = $wb.Worksheets.Item("control")
$wsControl.Range("B21").Value = 1
$wsControl.Range("B21").Select()
$wsControl.Activate()
